$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update event text cells (columns C and D) for rows 2-6
$ws.Range("C2").Value = "Milwaukee Bucks vs Indiana Pacers 2021-10-25 23:00:00+00:00Port-aux-Francais, Kerguelen,  TF.Category: sports"
$ws.Range("D2").Value = "PATINS Access to Education 2021-11-17 14:00:00+00:00Port-aux-Francais, Kerguelen,  TF.Category: conferences"

$ws.Range("C3").Value = "Passafire and Indubious 2021-11-08 01:00:00+00:00Port-aux-Francais, Kerguelen,  TF.Category: concerts"

$ws.Range("C4").Value = "State Ballet of Ukraine: Cinderella 2021-11-17 00:30:00+00:00Port-aux-Francais, Kerguelen,  TF.Category: performing-arts"

$ws.Range("B5").Value = "Major Restaurant"
$ws.Range("C5").Value = "The (W)hole in Our HeArts, part of Spirit & Place Festival 2021-11-06 14:00:00+00:00Port-aux-Francais, Kerguelen,  TF.Category: festivals"

$ws.Range("C6").Value = "Footnote & Risk Factor Disclosures: Current Examples & Best Practices 2021-11-17 14:00:00+00:00Port-aux-Francais, Kerguelen,  TF.Category: expos"

# Update restaurant rows (columns A and B) for rows 7-15
$ws.Range("B7").Value = "Bluebeard"
$ws.Range("B8").Value = "Yats"

$ws.Range("A9").Value = 12
$ws.Range("B9").Value = "Tinker Street Restaurant"

$ws.Range("A10").Value = 13
$ws.Range("B10").Value = "Livery"

$ws.Range("A11").Value = 18
$ws.Range("B11").Value = "Nesso"

$ws.Range("A12").Value = 19
$ws.Range("B12").Value = "Axum Ethiopian Restaurant"

$ws.Range("A13").Value = 20
$ws.Range("B13").Value = "Oakleys Bistro"

$ws.Range("A14").Value = 39
$ws.Range("B14").Value = "Mama Carolla's"

$ws.Range("A15").Value = 53
$ws.Range("B15").Value = "Kuma's Corner"

# Remove row 16 entirely (was Mama Carolla's / 59 before the shift)
$ws.Rows.Item(16).Delete()
